# Update speedup table workbook:
#  - add two new worksheets ("ijcnn1", "generated") after the existing "a1a" sheet
#  - populate them with the evaluation results
#  - update view/selection state to match the final authored state

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- create the two new sheets, in order, right after "a1a" ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "ijcnn1"

$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "generated"

# --- populate "ijcnn1" (sheet2) ---
$ws2.Range("B1").Value = "train time"
$ws2.Range("C1").Value = "Compress time in train"
$ws2.Range("D1").Value = "Split time in train"
$ws2.Range("E1").Value = "communication time"
$ws2.Range("F1").Value = "initCUDA time in train"
$ws2.Range("G1").Value = "test time"
$ws2.Range("H1").Value = "correct rate"

$ws2.Range("A2").Value = "Sequential"
$ws2.Range("B2").Value = 2.5315059999999998
$ws2.Range("C2").Value = 1.8905190000000001
$ws2.Range("D2").Value = 0.28919299999999998
$ws2.Range("G2").Value = 0.54247800000000002
$ws2.Range("H2").Value = 0.67373300000000003

$ws2.Range("A3").Value = "OpenMP - data parallel"
$ws2.Range("B3").Value = 2.4280279999999999
$ws2.Range("C3").Value = 1.563461
$ws2.Range("D3").Value = 0.53073899999999996
$ws2.Range("G3").Value = 0.65548200000000001
$ws2.Range("H3").Value = 0.67373300000000003

$ws2.Range("A4").Value = "OpenMP - feature parallel"

$ws2.Range("A5").Value = "OpenMPI"
$ws2.Range("B5").Value = 17.892775
$ws2.Range("C5").Value = 6.2186599999999999
$ws2.Range("D5").Value = 22.483117
$ws2.Range("E5").Value = [double]"1.1E-5"
$ws2.Range("H5").Value = 0.95622700000000005

$ws2.Range("A6").Value = "CUDA"
$ws2.Range("B6").Value = 3.6024859999999999
$ws2.Range("C6").Value = 3.1373679999999999
$ws2.Range("D6").Value = [double]"8.6230000000000005E-3"
$ws2.Range("F6").Value = 0.113554
$ws2.Range("G6").Value = 0.44085099999999999
$ws2.Range("H6").Value = 0.90499600000000002

# column E on "ijcnn1" was manually widened to fit "communication time"
$ws2.Columns.Item(5).ColumnWidth = 8.43

# a stray formatted row far below the data (row 22) carries a custom height
$ws2.Rows.Item(22).RowHeight = 14.5

# --- populate "generated" (sheet3) ---
$ws3.Range("B1").Value = "train time"
$ws3.Range("C1").Value = "Compress time in train"
$ws3.Range("D1").Value = "Split time in train"
$ws3.Range("E1").Value = "communication time"
$ws3.Range("F1").Value = "initCUDA time in train"
$ws3.Range("G1").Value = "test time"
$ws3.Range("H1").Value = "correct rate"

$ws3.Range("A2").Value = "Sequential"
$ws3.Range("B2").Value = 20.206942000000002
$ws3.Range("C2").Value = 13.738770000000001
$ws3.Range("D2").Value = 4.6531690000000001
$ws3.Range("G2").Value = 0.44645600000000002
$ws3.Range("H2").Value = 0.59909999999999997

$ws3.Range("A3").Value = "OpenMP - data parallel"
$ws3.Range("B3").Value = 13.650181
$ws3.Range("C3").Value = 6.5865970000000003
$ws3.Range("D3").Value = 5.2281789999999999
$ws3.Range("G3").Value = 0.47265299999999999
$ws3.Range("H3").Value = 0.59909999999999997

$ws3.Range("A4").Value = "OpenMP - feature parallel"

$ws3.Range("A5").Value = "OpenMPI"
$ws3.Range("B5").Value = 203.485626
$ws3.Range("C5").Value = 40.328243999999998
$ws3.Range("D5").Value = 322.14626700000002
# this particular value was entered as literal text (European-style decimal comma)
$ws3.Range("E5").NumberFormat = "@"
$ws3.Range("E5").Value = "0,000011"
$ws3.Range("E5").ClearFormats()
$ws3.Range("G5").Value = 0.61116899999999996
$ws3.Range("H5").Value = 0.77229999999999999

$ws3.Range("A6").Value = "CUDA"
$ws3.Range("B6").Value = 12.098867
$ws3.Range("C6").Value = 9.8292479999999998
$ws3.Range("D6").Value = 0.498996
$ws3.Range("F6").Value = 0.10288600000000001
$ws3.Range("G6").Value = 0.37675900000000001
$ws3.Range("H6").Value = 0.63790000000000002

# --- update selections per-sheet and make "generated" the active tab ---
$ws1.Activate()
$ws1.Range("B23").Select()

$ws2.Activate()
$ws2.Range("F6").Select()

$ws3.Activate()
$ws3.Range("H6").Select()
